# Update the "prediction" score column (B) on the active sheet
# ("quadratic-svm-score") with the refreshed quadratic-SVM distance
# values from the regenerated ful-path.csv output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1197.6748946188382
$ws.Range("B3").Value = 1204.6856609442207
$ws.Range("B4").Value = 1146.9521974307906
$ws.Range("B5").Value = 1344.2913157794121
$ws.Range("B6").Value = 1348.1430359773221
